$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'7.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.89%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.339"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'5.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07524"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'12.69%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.807"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'7.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'16.35%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9066"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.35%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01687"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,508.73%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1690"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.74%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07674"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'12.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08072"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.35%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03026"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.05%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09868"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'9.78%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001524"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.66%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04548"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.78%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006489"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.77%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.503"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.47%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.69%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'2.00%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.176"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.69%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'2.87%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001215"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.01%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004489"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'9.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'8.30%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'7.56%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04555"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'7.86%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007154"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.43%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1363"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'10.09%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'8.00%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01388"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.69%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006103"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.83%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-3.81%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.65%"
$ws.Range("E47").Style = "Normal"
